# Update column G ("K") values on the active sheet with newly regenerated
# strikeout (K) counts, replacing the previous "Strike#" derived values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(0, 8, 7, 6, 5, 6, 6, 2, 1, 6, 4, 5, 5, 2, 7, 7, 7, 3, 6, 7, 4, 3, 1, 10, 6, 7, 5, 5, 4, 5, 5)

for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
